$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (roll_offset stays, contract_cycle is new,
# old commision/multiplier shift right by one column).
$ws.Range("E1:E10").EntireColumn.Insert()

# Header
$ws.Range("E1").Value = "contract_cycle"

# Contract cycle values per symbol (row order matches existing sheet rows 2-10)
$ws.Range("E2").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E3").Value = "1,3,5,7,8,9,11,12"
$ws.Range("E4").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E5").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E6").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E7").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E8").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("E9").Value = "1,3,5,7,8,9,11"
$ws.Range("E10").Value = "1,2,3,4,5,6,7,8,9,10,11,12"

# Copy header style (bold, centered) to the new header cell to match other headers
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# Widen the new contract_cycle column (F/G already inherited the old E/F widths
# verbatim from the column insert above, so only E needs an explicit width).
$ws.Range("E1").EntireColumn.ColumnWidth = 27.15

$ws.Range("G21").Select()
